$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the literal text (otherwise values like "$8.54" get
    # auto-detected as currency numbers by Excel's input parser), then
    # drop the number-format override so no stray style sticks to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Updated On-Demand / SUD / 1-Year / 3-Year URL for row 2 (new dl= token)
$row2Url = "https://cloud.google.com/products/calculator?dl=CjhDaVJrWm1RMFlUZ3hOaTFtWlRNMkxUUXdORGN0T0dGbU9TMHpZelpoTkRKaVl6Y3labVVRQVE9PRAIGiRENUU4OUJFQi04NzhFLTQ3N0ItODlGMC0wMkQ2RjkxMjdGREI"

$ws.Range("E2").Value = $row2Url
$ws.Range("G2").Value = $row2Url
$ws.Range("I2").Value = $row2Url
$ws.Range("K2").Value = $row2Url

# Row 2 was missing its 3-Year Price (L2) -- fill it in
Set-TextValue $ws.Range("L2") "`$8.54"

# New row 3
$row3Url = "https://cloud.google.com/products/calculator?dl=CjhDaVJsTXpkbVpUazROaTFpTldOa0xUUXhZek10T0dJd09TMWtNVGxoTVdVMFlqRTBPV1VRQVE9PRAIGiQ2RkNCRThFRi03OUI4LTRFOUYtODcxRS1EODg3NjhEMTc3QkU"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Free: Debian, CentOS, CoreOS, Ubuntu or BY"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "general purpose"
$ws.Range("E3").Value = $row3Url
Set-TextValue $ws.Range("F3") "`$8.54"
$ws.Range("G3").Value = $row3Url
Set-TextValue $ws.Range("H3") "`$8.54"
$ws.Range("I3").Value = $row3Url
Set-TextValue $ws.Range("J3") "`$8.54"
$ws.Range("K3").Value = $row3Url
Set-TextValue $ws.Range("L3") "`$8.54"

# New row 4
$row4Url = "https://cloud.google.com/products/calculator?dl=CjhDaVF3WXpCbE9EZGlNUzB4WTJFeUxUUTROVEl0WW1VMk9DMWtZVE16TWpReE9HVTVZalFRQVE9PRAIGiRFMUY2OTVEQy03NjExLTQ3MzktOTBFNC0zMzVDQjIyMzk3NTc"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Free: Debian, CentOS, CoreOS, Ubuntu or BY"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "general purpose"
$ws.Range("E4").Value = $row4Url
Set-TextValue $ws.Range("F4") "`$8.54"
$ws.Range("G4").Value = $row4Url
Set-TextValue $ws.Range("H4") "`$8.54"
$ws.Range("I4").Value = $row4Url
Set-TextValue $ws.Range("J4") "`$8.54"
$ws.Range("K4").Value = $row4Url
Set-TextValue $ws.Range("L4") "`$8.54"
